# Insert a new weekly record row at position 94 on the single sheet,
# shifting the existing rows 94-116 down to 95-117.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(94).Insert()

# Populate the newly inserted row 94 with the new weekly observation.
$ws.Cells.Item(94, 1).Value = 4
$ws.Cells.Item(94, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(94, 3).Value = "Los Lagos"
$ws.Cells.Item(94, 4).Value = 44932
$ws.Cells.Item(94, 5).Value = 10
$ws.Cells.Item(94, 6).Value = 100112031
$ws.Cells.Item(94, 7).Value = "Poroto verde"
$ws.Cells.Item(94, 8).Value = "Magnum"
$ws.Cells.Item(94, 9).Value = "Primera"
$ws.Cells.Item(94, 10).Value = 35
$ws.Cells.Item(94, 11).Value = 35000
$ws.Cells.Item(94, 12).Value = 35000
$ws.Cells.Item(94, 13).Value = 35000
$ws.Cells.Item(94, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(94, 15).Value = "Región Metropolitana"
$ws.Cells.Item(94, 16).Value = 1400
$ws.Cells.Item(94, 17).Value = 25
$ws.Cells.Item(94, 18).Value = "Hortaliza"

Write-Host "Row 94 inserted and populated"
